$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update values reflecting the new test scenarios
$ws.Range("D2").Value = "Automation5"
$ws.Range("F2").Value = "Sanity.xlsx"

# Widen column D (ColumnWidth in characters maps to stored width=13 units)
$ws.Columns.Item(4).ColumnWidth = 12.142857142857142

# Move active selection to G2
$ws.Range("G2").Select()
